$d = $word.ActiveDocument

# Locate the target run: "summer internship that utilizes ... through computing."
$targetText = "summer internship that utilizes my computational, analytical, communicative, and teamwork skills to solve problems and make a positive difference in my community through computing."
$content = $d.Content.Text
$startIdx = $content.IndexOf($targetText)

$fullRange = $d.Range($startIdx, $startIdx + $targetText.Length)

# Shrink the existing run down to just the trailing portion first so the run
# is never left empty (avoids corrupting later range/bookmark resolution).
$fullRange.Text = "that utilizes my computational, analytical, communicative, and teamwork skills to solve problems and make a positive difference in my community through computing."

# Place the (re-homed) _GoBack bookmark right before "that utilizes...".
# Adding a bookmark with a name that already exists elsewhere in the
# document moves it here instead of creating a duplicate.
$bmPoint = $d.Range($startIdx, $startIdx)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# Insert "internship " immediately before the bookmark/"that utilizes..." text;
# this also carries the bookmark along with it to stay right after "internship ".
$insertPoint = $d.Range($startIdx, $startIdx)
$insertPoint.InsertBefore("internship ")

# Insert "summer " before "internship ".
$content2 = $d.Content.Text
$internshipIdx = $content2.IndexOf("internship ", $startIdx - 1)
$insertPoint2 = $d.Range($internshipIdx, $internshipIdx)
$insertPoint2.InsertBefore("summer ")

# Force "summer " to stay a distinct run (rather than being silently
# re-merged with its identically formatted neighbours on save) by toggling
# a character property on and back off.
$content3 = $d.Content.Text
$summerIdx = $content3.IndexOf("summer ", $startIdx - 1)
$summerRange = $d.Range($summerIdx, $summerIdx + "summer ".Length)
$summerRange.Bold = 1
$summerRange.Bold = 0
